# Update "想去人数" (interest count) figures across the refreshed data pull.
$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheetId 1)
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F3").Value = 234
$wsExpo.Range("F4").Value = 820
$wsExpo.Range("F6").Value = 405
$wsExpo.Range("F7").Value = 567
$wsExpo.Range("F11").Value = 131
$wsExpo.Range("F14").Value = 1767
$wsExpo.Range("F15").Value = 325
$wsExpo.Range("F16").Value = 2523
$wsExpo.Range("F17").Value = 301

# Sheet "本地生活" (sheetId 3)
$wsLocal = $wb.Worksheets.Item("本地生活")
$wsLocal.Range("F2").Value = 5292
$wsLocal.Range("F4").Value = 207

# Sheet "全部类型" (sheetId 4)
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F3").Value = 5292
$wsAll.Range("F6").Value = 207
$wsAll.Range("F7").Value = 234
$wsAll.Range("F13").Value = 820
$wsAll.Range("F17").Value = 405
$wsAll.Range("F18").Value = 567
$wsAll.Range("F23").Value = 131
$wsAll.Range("F29").Value = 1767
$wsAll.Range("F30").Value = 325
$wsAll.Range("F31").Value = 2524
$wsAll.Range("F33").Value = 301
